# Fruta / hortaliza, semanal
# Updates the Achicoria price records (rows 3-14) to reflect the new weekly
# data pull: Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M), Origen (O) and Precio $/Kg (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44230
$ws.Range("J3").Value = 250

$ws.Range("D4").Value = 44251
$ws.Range("J4").Value = 120
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 5000
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 312

$ws.Range("D5").Value = 44186
$ws.Range("J5").Value = 160

$ws.Range("D6").Value = 44188
$ws.Range("J6").Value = 210

$ws.Range("D7").Value = 44232
$ws.Range("J7").Value = 250

$ws.Range("D8").Value = 44189
$ws.Range("J8").Value = 250

$ws.Range("D9").Value = 44215
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 5000
$ws.Range("M9").Value = 5500
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 344

$ws.Range("D10").Value = 44210
$ws.Range("J10").Value = 340

$ws.Range("D11").Value = 44231
$ws.Range("J11").Value = 250
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 5500
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 344

$ws.Range("D12").Value = 44204
$ws.Range("J12").Value = 430

$ws.Range("D13").Value = 44187
$ws.Range("J13").Value = 160

$ws.Range("D14").Value = 44292
$ws.Range("J14").Value = 90
$ws.Range("K14").Value = 6000
$ws.Range("M14").Value = 6000
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 375
